$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "heats": extend the data series with one more point (column Q) and
# insert a new leading "0" dilution/time-zero column for the observation and
# dilution rows (shifting their existing values one column to the right).
# ---------------------------------------------------------------------------
$heats = $wb.Worksheets.Item("heats")

# Row 1 ("data" / point index) and Row 2 ("volumes") simply gain one more
# trailing value.
$heats.Range("Q1").Value = 16
$heats.Range("Q2").Value = 15.374955

# Row 3 ("observation"): shift existing B3:P3 -> C3:Q3, then set the new
# leading value to 0.
$heats.Range("Q3").Value = $heats.Range("P3").Value()
$heats.Range("P3").Value = $heats.Range("O3").Value()
$heats.Range("O3").Value = $heats.Range("N3").Value()
$heats.Range("N3").Value = $heats.Range("M3").Value()
$heats.Range("M3").Value = $heats.Range("L3").Value()
$heats.Range("L3").Value = $heats.Range("K3").Value()
$heats.Range("K3").Value = $heats.Range("J3").Value()
$heats.Range("J3").Value = $heats.Range("I3").Value()
$heats.Range("I3").Value = $heats.Range("H3").Value()
$heats.Range("H3").Value = $heats.Range("G3").Value()
$heats.Range("G3").Value = $heats.Range("F3").Value()
$heats.Range("F3").Value = $heats.Range("E3").Value()
$heats.Range("E3").Value = $heats.Range("D3").Value()
$heats.Range("D3").Value = $heats.Range("C3").Value()
$heats.Range("C3").Value = $heats.Range("B3").Value()
$heats.Range("B3").Value = 0

# Row 4 ("dilution"): shift existing B4:P4 -> C4:Q4, then set the new
# leading value to 0.
$heats.Range("Q4").Value = $heats.Range("P4").Value()
$heats.Range("P4").Value = $heats.Range("O4").Value()
$heats.Range("O4").Value = $heats.Range("N4").Value()
$heats.Range("N4").Value = $heats.Range("M4").Value()
$heats.Range("M4").Value = $heats.Range("L4").Value()
$heats.Range("L4").Value = $heats.Range("K4").Value()
$heats.Range("K4").Value = $heats.Range("J4").Value()
$heats.Range("J4").Value = $heats.Range("I4").Value()
$heats.Range("I4").Value = $heats.Range("H4").Value()
$heats.Range("H4").Value = $heats.Range("G4").Value()
$heats.Range("G4").Value = $heats.Range("F4").Value()
$heats.Range("F4").Value = $heats.Range("E4").Value()
$heats.Range("E4").Value = $heats.Range("D4").Value()
$heats.Range("D4").Value = $heats.Range("C4").Value()
$heats.Range("C4").Value = $heats.Range("B4").Value()
$heats.Range("B4").Value = 0

# Row 5 ("deviation"): constant value series, simply gains a trailing value
# (same constant as its neighbours).
$heats.Range("Q5").Value = $heats.Range("P5").Value()

# Tidy up the selection on this sheet now that it has an extra column.
$heats.Range("R7").Select()

# ---------------------------------------------------------------------------
# Sheet "enthalpies": the extra "Comp / 47200" row is no longer needed.
# ---------------------------------------------------------------------------
$enthalpies = $wb.Worksheets.Item("enthalpies")
$enthalpies.Rows("2:2").Delete()
$enthalpies.Range("A2").Select()

# ---------------------------------------------------------------------------
# Sheet "input_concentrations": tidy the selection (no data changes).
# ---------------------------------------------------------------------------
$concentrations = $wb.Worksheets.Item("input_concentrations")
$concentrations.Range("A3").Select()

# The "enthalpies" sheet is now the active/visible tab (was "heats" before).
$enthalpies.Activate()
